$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Updated Training Set data - fill in rows 77-89 (columns A, B, C)
# Column A gets the new temperature/pressure values; B and C are zeroed out
$values = @(
    @(373, 0, 0),
    @(383, 0, 0),
    @(393, 0, 0),
    @(403, 0, 0),
    @(413, 0, 0),
    @(423, 0, 0),
    @(433, 0, 0),
    @(443, 0, 0),
    @(453, 0, 0),
    @(463, 0, 0),
    @(473, 0, 0),
    @(483, 0, 0),
    @(493, 0, 0)
)

$startRow = 77
for ($i = 0; $i -lt $values.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $values[$i][0]
    $ws.Cells.Item($row, 2).Value = $values[$i][1]
    $ws.Cells.Item($row, 3).Value = $values[$i][2]
}

# Update the active sheet's selection / scroll position to match where the
# author was working when the training data was added
$ws.Activate()
$ws.Range("K72").Select()
